# Regenerate the handback report timestamps (commit: "Generate Report for Handback")
$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-13 09:07:32"
$wsZhCn.Range("H2").Value = "2016-03-13 09:07:49"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-13 09:07:36"
$wsDeDe.Range("H2").Value = "2016-03-13 09:07:55"
